$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Column A "coupon name" cells lose their generic "كوبون"/"Coupon" wrapper ---
# Arabic block (rows 26-39)
$ws.Range("A26").Value = "بلندز"
$ws.Range("A27").Value = "تيمو"
$ws.Range("A28").Value = "كوهوم"
$ws.Range("A29").Value = "بيتونيا"
$ws.Range("A30").Value = "جلوبال يو"
$ws.Range("A31").Value = "ون كارد"
$ws.Range("A32").Value = "بازل إنجلش"
$ws.Range("A33").Value = "كوبون سيفي"
$ws.Range("A34").Value = "مجموعة طبيب"
$ws.Range("A35").Value = "نسبة"
$ws.Range("A36").Value = "ستور اص"
$ws.Range("A37").Value = "سويس اربيان"
$ws.Range("A38").Value = "فسرلي"
$ws.Range("A39").Value = "بينه"

# English block (rows 83-96)
$ws.Range("A83").Value = "Blends"
$ws.Range("A84").Value = "Temu"
$ws.Range("A85").Value = "Coohom"
$ws.Range("A86").Value = "Baytonia"
$ws.Range("A87").Value = "Global YO"
$ws.Range("A88").Value = "OneCard"
$ws.Range("A89").Value = "Puzzle English"
$ws.Range("A90").Value = "Sivvi"
$ws.Range("A91").Value = "Tabib Group"
$ws.Range("A92").Value = "Nesbh"
$ws.Range("A93").Value = "Storeus"
$ws.Range("A94").Value = "Swiss Arabian"
$ws.Range("A95").Value = "Fasrly"
$ws.Range("A96").Value = "Bynh"

# --- Scroll/zoom/selection state of the view ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 230
$ws.Range("B99").Select()
